$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subject names (column G) to fill in for rows 3-16
$subjects = @{
    3  = "Intro"
    4  = "Linear Algebra"
    5  = "Statistics"
    6  = "ML"
    7  = "ML"
    8  = "Linear Regression"
    9  = "Gradient Descent"
    10 = "Polynominal Regression"
    11 = "Over and under fitting"
    12 = "Hold out and crossvalidation"
    13 = "Regularization"
    14 = "Logistical regression"
    15 = "linear classification"
    16 = "FLDA"
}

for ($r = 3; $r -le 17; $r++) {
    $ws.Range("D$r").Value = "-"
    $ws.Range("E$r").Value = "-"
    if ($subjects.ContainsKey($r)) {
        $ws.Range("G$r").Value = $subjects[$r]
    }
}

# Freeze panes: split at column D / row 9 (3 columns, 8 rows frozen),
# with the active selection in the bottom-right pane at G17.
$ws.Activate()
$ws.Range("D9").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("G17").Select()
